# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Headers go in AD1:AF1 (matching the formatting of the existing header row),
# and every data row (2-35) gets the team's season record: 53 wins, 60 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Give the new header cells the same formatting as the rest of row 1
# (bold font, centered, thin border) by copying the format from AC1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the season record for every data row
for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 30).Value = 53  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 60  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-35"
